$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 9 block ---
# D9 and F9 keep the same displayed text ("1.107/1.12" and "1.076/1.1000"),
# only their underlying shared-string index shifted because two strings
# were removed elsewhere. G9's numeric flag is cleared.
$ws.Range("G9").Value = $null

# --- Row 18 block ---
# B18 text changes from "5.005/5.029" to "4.988/5.012".
# C18, D18 and F18 are cleared (their old values/strings were removed).
$ws.Range("B18").Value = "4.988/5.012"
$ws.Range("C18").Value = $null
$ws.Range("D18").Value = $null
$ws.Range("F18").Value = $null

# --- Row 21 block ---
# B21 keeps the same displayed text ("科创50（588000）"); only the
# shared-string index shifted.

# --- Row 26 block ---
# B26 text changes from "1.412/1.444" to "1.416/1.444".
$ws.Range("B26").Value = "1.416/1.444"

# --- Row 27 block ---
# B27 and D27 keep the same displayed text ("1.452/1.464" and
# "1.430/1.445"); only the shared-string index shifted.

# --- View state ---
# Scroll the visible window so row 10 is the top-left visible row (was row 19)
# and move/select the active cell to D12 (was D32).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
